$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.78870267507108
$ws.Range("C2").Value = 4.739763129625044
$ws.Range("D2").Value = 10.39049875215931
$ws.Range("F2").Value = 34.72366958286326
$ws.Range("G2").Value = 3.664517566322587
$ws.Range("J2").Value = 11.40098933484022
$ws.Range("K2").Value = 10.11336487028463
$ws.Range("O2").Value = 25.76425711858706
$ws.Range("B3").Value = 10.50461739567848
$ws.Range("C3").Value = 4.53753218069159
$ws.Range("D3").Value = 10.28506853920411
$ws.Range("F3").Value = 34.77864920164437
$ws.Range("G3").Value = 3.66649278908222
$ws.Range("J3").Value = 11.38327850150298
$ws.Range("K3").Value = 9.914647745785389
$ws.Range("O3").Value = 25.85619558417756
$ws.Range("B4").Value = 10.32802025247671
$ws.Range("C4").Value = 4.407875623387196
$ws.Range("D4").Value = 10.22186733867488
$ws.Range("F4").Value = 34.82112680498717
$ws.Range("G4").Value = 3.667769451852884
$ws.Range("J4").Value = 11.37468163205411
$ws.Range("K4").Value = 9.792046724735156
$ws.Range("O4").Value = 25.91820972920383
$ws.Range("B5").Value = 10.25562011887286
$ws.Range("C5").Value = 4.35370790196287
$ws.Range("D5").Value = 10.19652431176859
$ws.Range("F5").Value = 34.84062428721671
$ws.Range("G5").Value = 3.66830581374136
$ws.Range("J5").Value = 11.37175354997994
$ws.Range("K5").Value = 9.7420040413874
$ws.Range("O5").Value = 25.94487679078379
$ws.Range("B6").Value = 10.24357515394063
$ws.Range("C6").Value = 4.344634504536236
$ws.Range("D6").Value = 10.19234175951974
$ws.Range("F6").Value = 34.8439937829292
$ws.Range("G6").Value = 3.668395850825021
$ws.Range("J6").Value = 11.37130214938952
$ws.Range("K6").Value = 9.733691515218199
$ws.Range("O6").Value = 25.94938904352141
$ws.Range("B7").Value = 10.32704545259299
$ws.Range("C7").Value = 4.40715042181025
$ws.Range("D7").Value = 10.22152385139128
$ws.Range("F7").Value = 34.8213809051602
$ws.Range("G7").Value = 3.667776620114287
$ws.Range("J7").Value = 11.37463981095529
$ws.Range("K7").Value = 9.791372074226496
$ws.Range("O7").Value = 25.91856372332444
$ws.Range("B8").Value = 10.69126397472929
$ws.Range("C8").Value = 4.671199087203686
$ws.Range("D8").Value = 10.35384422791574
$ws.Range("F8").Value = 34.74081386193909
$ws.Range("G8").Value = 3.66518539826976
$ws.Range("J8").Value = 11.39441152621359
$ws.Range("K8").Value = 10.04500832207051
$ws.Range("O8").Value = 25.79480087847364
$ws.Range("B9").Value = 11.38378292264876
$ws.Range("C9").Value = 5.143703677744011
$ws.Range("D9").Value = 10.62429779148237
$ws.Range("F9").Value = 34.65220951311531
$ws.Range("G9").Value = 3.660608439839609
$ws.Range("J9").Value = 11.45112913040344
$ws.Range("K9").Value = 10.53491718458097
$ws.Range("O9").Value = 25.59638760211478
$ws.Range("B10").Value = 11.8737595720965
$ws.Range("C10").Value = 5.461290053835364
$ws.Range("D10").Value = 10.82805774013231
$ws.Range("F10").Value = 34.62963566491945
$ws.Range("G10").Value = 3.65754997075908
$ws.Range("J10").Value = 11.50354352218873
$ws.Range("K10").Value = 10.88670575662131
$ws.Range("O10").Value = 25.47779998042799
$ws.Range("B11").Value = 12.09153472538518
$ws.Range("C11").Value = 5.599038336677794
$ws.Range("D11").Value = 10.92149947420113
$ws.Range("F11").Value = 34.62862924842355
$ws.Range("G11").Value = 3.656223953702722
$ws.Range("J11").Value = 11.52967040918007
$ws.Range("K11").Value = 11.04426893978474
$ws.Range("O11").Value = 25.42979371650409
$ws.Range("B12").Value = 12.17318726570586
$ws.Range("C12").Value = 5.650212739735719
$ws.Range("D12").Value = 10.95696180377295
$ws.Range("F12").Value = 34.62958082157746
$ws.Range("G12").Value = 3.655731162052665
$ws.Range("J12").Value = 11.53988731373015
$ws.Range("K12").Value = 11.10352463987123
$ws.Range("O12").Value = 25.41247223896173
$ws.Range("B13").Value = 12.15563931374295
$ws.Range("C13").Value = 5.639235650588665
$ws.Range("D13").Value = 10.94932136231807
$ws.Range("F13").Value = 34.62931661073424
$ws.Range("G13").Value = 3.655836878801187
$ws.Range("J13").Value = 11.53767263115725
$ws.Range("K13").Value = 11.09078194639324
$ws.Range("O13").Value = 25.41616454203783
$ws.Range("B14").Value = 12.09826907011781
$ws.Range("C14").Value = 5.603268392297322
$ws.Range("D14").Value = 10.92441558916102
$ws.Range("F14").Value = 34.62868082530745
$ws.Range("G14").Value = 3.656183224473842
$ws.Range("J14").Value = 11.53050451305525
$ws.Range("K14").Value = 11.04915246471911
$ws.Range("O14").Value = 25.42835146490621
$ws.Range("B15").Value = 12.06301990429422
$ws.Range("C15").Value = 5.581108192302657
$ws.Range("D15").Value = 10.90916935465417
$ws.Range("F15").Value = 34.62846494777138
$ws.Range("G15").Value = 3.656396586437037
$ws.Range("J15").Value = 11.52615576747349
$ws.Range("K15").Value = 11.02359822037311
$ws.Range("O15").Value = 25.43592806916276
$ws.Range("B16").Value = 11.85941762319838
$ws.Range("C16").Value = 5.452150704213592
$ws.Range("D16").Value = 10.82196363821864
$ws.Range("F16").Value = 34.62988790844493
$ws.Range("G16").Value = 3.657637938965715
$ws.Range("J16").Value = 11.50188158680609
$ws.Range("K16").Value = 10.87635418394549
$ws.Range("O16").Value = 25.48105711286712
$ws.Range("B17").Value = 11.73314555223447
$ws.Range("C17").Value = 5.371301286540469
$ws.Range("D17").Value = 10.76863657218553
$ws.Range("F17").Value = 34.63313404796727
$ws.Range("G17").Value = 3.658416158648331
$ws.Range("J17").Value = 11.48757146636977
$ws.Range("K17").Value = 10.78535243856224
$ws.Range("O17").Value = 25.51026604383523
$ws.Range("B18").Value = 11.66003972753972
$ws.Range("C18").Value = 5.324167016508008
$ws.Range("D18").Value = 10.73803665828647
$ws.Range("F18").Value = 34.6358729938168
$ws.Range("G18").Value = 3.658869918783619
$ws.Range("J18").Value = 11.47955591346514
$ws.Range("K18").Value = 10.73278151380003
$ws.Range("O18").Value = 25.52762514253485
$ws.Range("B19").Value = 11.63520784477279
$ws.Range("C19").Value = 5.308100348328781
$ws.Range("D19").Value = 10.72768940399735
$ws.Range("F19").Value = 34.63695005657555
$ws.Range("G19").Value = 3.659024611668822
$ws.Range("J19").Value = 11.47687910347823
$ws.Range("K19").Value = 10.71494428810868
$ws.Range("O19").Value = 25.53359852314173
$ws.Range("B20").Value = 11.7466374350595
$ws.Range("C20").Value = 5.379973400505545
$ws.Range("D20").Value = 10.77430605335519
$ws.Range("F20").Value = 34.63269824889211
$ws.Range("G20").Value = 3.658332679804995
$ws.Range("J20").Value = 11.48907256132539
$ws.Range("K20").Value = 10.79506385461601
$ws.Range("O20").Value = 25.50709883547306
$ws.Range("B21").Value = 12.11514278740931
$ws.Range("C21").Value = 5.613859809045278
$ws.Range("D21").Value = 10.93172913946202
$ws.Range("F21").Value = 34.62883140135411
$ws.Range("G21").Value = 3.656081241183468
$ws.Range("J21").Value = 11.5326012338601
$ws.Range("K21").Value = 11.06139160529345
$ws.Range("O21").Value = 25.42474857336226
$ws.Range("B22").Value = 12.35120519803873
$ws.Range("C22").Value = 5.760952606065094
$ws.Range("D22").Value = 11.03505504980555
$ws.Range("F22").Value = 34.63407171627269
$ws.Range("G22").Value = 3.654664227903
$ws.Range("J22").Value = 11.56293104529072
$ws.Range("K22").Value = 11.2330388051663
$ws.Range("O22").Value = 25.37592726691269
$ws.Range("B23").Value = 12.22567542925733
$ws.Range("C23").Value = 5.682980017423545
$ws.Range("D23").Value = 10.97987754163154
$ws.Range("F23").Value = 34.63056414638005
$ws.Range("G23").Value = 3.655415549539491
$ws.Range("J23").Value = 11.5465731114563
$ws.Range("K23").Value = 11.14166567747055
$ws.Range("O23").Value = 25.40152561004825
$ws.Range("B24").Value = 11.74053933837007
$ws.Range("C24").Value = 5.376054769702169
$ws.Range("D24").Value = 10.77174269633899
$ws.Range("F24").Value = 34.63289255528135
$ws.Range("G24").Value = 3.658370400802271
$ws.Range("J24").Value = 11.48839325703596
$ws.Range("K24").Value = 10.79067410661705
$ws.Range("O24").Value = 25.50852896598198
$ws.Range("B25").Value = 11.19936552742723
$ws.Range("C25").Value = 5.020949698838963
$ws.Range("D25").Value = 10.55012533465177
$ws.Range("F25").Value = 34.668724325157
$ws.Range("G25").Value = 3.661792967731498
$ws.Range("J25").Value = 11.43388282434614
$ws.Range("K25").Value = 10.40356502608355
$ws.Range("O25").Value = 25.64530275865788
